$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "sahilharpal1234@gmail.com"
$ws.Range("B3").Value = "Sahil Harpal"
$ws.Range("C3").Value = "CA245368"
$ws.Range("D3").Value = "7276801998"
$ws.Range("E3").Value = "Indian Institute of Technology Jodhpur"
